$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The F2 cell held a broken, unparsable formula (if/else Python-ish syntax
# isn't valid Excel). Replace it with the plain text "e" and make sure it
# picks up the same "general" cell style used by the rest of the sheet
# (rather than keeping its own one-off style).
$ws.Range("F2").Value = "e"
$ws.Range("F2").HorizontalAlignment = 1

# Column F no longer needs to be wide enough to fit the old formula text;
# resize it down to fit the remaining (much shorter) column contents.
$ws.Columns.Item(6).ColumnWidth = 24.3
